$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B3").Value = 29418
$ws.Range("K7").Select()
